$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new registered users, following the same layout as existing rows:
# Column A: documento (numeric id)
# Column B: Correo (email string)
# Column C: Celular (numeric id, same as documento)
# Column D: Contraseña (string, same placeholder password used by other rows)

$ws.Range("A53").Value = 66000128
$ws.Range("B53").Value = "66000128test@gmail.com"
$ws.Range("C53").Value = 66000128
$ws.Range("D53").Value = "Aaaaaaaaa1"

$ws.Range("A54").Value = 66000129
$ws.Range("B54").Value = "66000129test@gmail.com"
$ws.Range("C54").Value = 66000129
$ws.Range("D54").Value = "Aaaaaaaaa1"
